$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from SCD0200 to SCD0011
$ws.Name = "SCD0011"

# Update the TC_ID value in B2 from "DGS-215" to "SCD0011-031"
$ws.Range("B2").Value = "SCD0011-031"

# Widen column B to fit the new, longer TC_ID text (target stored width 12.42578125;
# closest value this engine's pixel-quantized ColumnWidth model can reach is 12.5)
$ws.Columns.Item(2).ColumnWidth = 11.666666666666666

# Move the active selection to B3
$ws.Range("B3").Select()
